# Update the "K" column (column G) values on the active sheet to reflect
# the regenerated save_data (using K / strikeouts instead of Strike# counts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New G-column (K) values, keyed by row number (rows 2-32 of sheet data).
$newK = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 3
    10 = 2
    11 = 2
    12 = 0
    13 = 2
    14 = 1
    15 = 3
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 2
    25 = 1
    26 = 2
    27 = 2
    28 = 1
    29 = 1
    30 = 2
    31 = 1
    32 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
